$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035269238710401
$ws.Range("D2").Value = 1.034677586403836
$ws.Range("E2").Value = 1.043146302407664
$ws.Range("F2").Value = 1.051224250301505
$ws.Range("I2").Value = 1.037381170959519
$ws.Range("J2").Value = 1.040383936404246
$ws.Range("K2").Value = 1.037476372845884
$ws.Range("L2").Value = 1.045920991050894
$ws.Range("M2").Value = 1.053976345408485
$ws.Range("N2").Value = 1.041861400927646
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036236780058784
$ws.Range("D3").Value = 1.035173233315323
$ws.Range("E3").Value = 1.044041305763808
$ws.Range("F3").Value = 1.052284352942675
$ws.Range("I3").Value = 1.037570943987283
$ws.Range("J3").Value = 1.040994864668283
$ws.Range("K3").Value = 1.037782092557411
$ws.Range("L3").Value = 1.0466267502221
$ws.Range("M3").Value = 1.05484841483963
$ws.Range("N3").Value = 1.042473196779892
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036863063637582
$ws.Range("D4").Value = 1.03549384550145
$ws.Range("E4").Value = 1.04462106386933
$ws.Range("F4").Value = 1.052971268778429
$ws.Range("I4").Value = 1.037692301514832
$ws.Range("J4").Value = 1.041389790432152
$ws.Range("K4").Value = 1.037979090886472
$ws.Range("L4").Value = 1.047083420515558
$ws.Range("M4").Value = 1.055413059150803
$ws.Range("N4").Value = 1.042868683383644
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037126405176462
$ws.Range("D5").Value = 1.035628604202486
$ws.Range("E5").Value = 1.044864944249649
$ws.Range("F5").Value = 1.053260276777452
$ws.Range("I5").Value = 1.037742975218946
$ws.Range("J5").Value = 1.04155572413521
$ws.Range("K5").Value = 1.038061710950928
$ws.Range("L5").Value = 1.047275403162099
$ws.Range("M5").Value = 1.0556505202617
$ws.Range("N5").Value = 1.043034852731593
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037170624405311
$ws.Range("D6").Value = 1.035651229169971
$ws.Range("E6").Value = 1.044905901609597
$ws.Range("F6").Value = 1.053308815870369
$ws.Range("I6").Value = 1.037751463306461
$ws.Range("J6").Value = 1.041583579653369
$ws.Range("K6").Value = 1.038075571579605
$ws.Range("L6").Value = 1.047307637779187
$ws.Range("M6").Value = 1.055690395966128
$ws.Range("N6").Value = 1.043062747807783
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03686658221691
$ws.Range("D7").Value = 1.035495646260653
$ws.Range("E7").Value = 1.044624322022492
$ws.Range("F7").Value = 1.052975129620404
$ws.Range("I7").Value = 1.0376929799758
$ws.Range("J7").Value = 1.041392008013136
$ws.Range("K7").Value = 1.037980195639292
$ws.Range("L7").Value = 1.04708598580416
$ws.Range("M7").Value = 1.055416231786256
$ws.Range("N7").Value = 1.042870904113848
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03559617796787
$ws.Range("D8").Value = 1.03484511325225
$ws.Range("E8").Value = 1.043448641703978
$ws.Range("F8").Value = 1.051582317942078
$ws.Range("I8").Value = 1.037445603054817
$ws.Range("J8").Value = 1.040590481873114
$ws.Range("K8").Value = 1.037579861976677
$ws.Range("L8").Value = 1.046159505820571
$ws.Range("M8").Value = 1.054270990914359
$ws.Range("N8").Value = 1.04206823971477
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.0333592742507
$ws.Range("D9").Value = 1.033698062832097
$ws.Range("E9").Value = 1.041381820593618
$ws.Range("F9").Value = 1.04913538641577
$ws.Range("I9").Value = 1.036998706838267
$ws.Range("J9").Value = 1.039175169826446
$ws.Range("K9").Value = 1.036868168155572
$ws.Range("L9").Value = 1.044526941552829
$ws.Range("M9").Value = 1.05225570195926
$ws.Range("N9").Value = 1.040650917762658
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03186919221322
$ws.Range("D10").Value = 1.032932970136607
$ws.Range("E10").Value = 1.040007277981717
$ws.Range("F10").Value = 1.04750911232073
$ws.Range("I10").Value = 1.036693428272533
$ws.Range("J10").Value = 1.038229708403607
$ws.Range("K10").Value = 1.036389562756816
$ws.Range("L10").Value = 1.043438619857763
$ws.Range("M10").Value = 1.050914087956599
$ws.Range("N10").Value = 1.039704113676152
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031224258910093
$ws.Range("D11").Value = 1.03260160388946
$ws.Range("E11").Value = 1.039412889523001
$ws.Range("F11").Value = 1.046806117092416
$ws.Range("I11").Value = 1.036559504518855
$ws.Range("J11").Value = 1.037819866545544
$ws.Range("K11").Value = 1.03618135175541
$ws.Range("L11").Value = 1.042967386869891
$ws.Range("M11").Value = 1.050333617481635
$ws.Range("N11").Value = 1.039293689795641
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030984744831956
$ws.Range("D12").Value = 1.03247850992409
$ws.Range("E12").Value = 1.039192227919833
$ws.Range("F12").Value = 1.046545173240338
$ws.Range("I12").Value = 1.036509499050869
$ws.Range("J12").Value = 1.037667565954589
$ws.Range("K12").Value = 1.036103867961148
$ws.Range("L12").Value = 1.042792353289391
$ws.Range("M12").Value = 1.050118074320074
$ws.Range("N12").Value = 1.039141172920378
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031036119467061
$ws.Range("D13").Value = 1.032504914414156
$ws.Range("E13").Value = 1.039239555104606
$ws.Range("F13").Value = 1.046601138430321
$ws.Range("I13").Value = 1.036520237159363
$ws.Range("J13").Value = 1.0377002379732
$ws.Range("K13").Value = 1.036120495041197
$ws.Range("L13").Value = 1.042829898424425
$ws.Range("M13").Value = 1.050164305906962
$ws.Range("N13").Value = 1.039173891337002
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031204459724365
$ws.Range("D14").Value = 1.032591429094083
$ws.Range("E14").Value = 1.039394647097811
$ws.Range("F14").Value = 1.046784543718346
$ws.Range("I14").Value = 1.03655537635622
$ws.Range("J14").Value = 1.037807278695744
$ws.Range("K14").Value = 1.036174949871959
$ws.Range("L14").Value = 1.042952918467454
$ws.Range("M14").Value = 1.050315799196363
$ws.Range("N14").Value = 1.03928108406965
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031308185407768
$ws.Range("D15").Value = 1.032644732395757
$ws.Range("E15").Value = 1.039490220419385
$ws.Range("F15").Value = 1.046897569627724
$ws.Range("I15").Value = 1.036576992311293
$ws.Range("J15").Value = 1.037873221148961
$ws.Range("K15").Value = 1.036208482112849
$ws.Range("L15").Value = 1.043028715638885
$ws.Range("M15").Value = 1.050409148432122
$ws.Range("N15").Value = 1.039347120168715
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031912000296884
$ws.Range("D16").Value = 1.032954960360786
$ws.Range("E16").Value = 1.040046742467621
$ws.Range("F16").Value = 1.047555792981942
$ws.Range("I16").Value = 1.036702279824998
$ws.Range("J16").Value = 1.038256898815481
$ws.Range("K16").Value = 1.036403360650564
$ws.Range("L16").Value = 1.043469894467859
$ws.Range("M16").Value = 1.050952621565516
$ws.Range("N16").Value = 1.03973134270153
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032290833283576
$ws.Range("D17").Value = 1.033149538964031
$ws.Range("E17").Value = 1.040396048218589
$ws.Range("F17").Value = 1.047968998823023
$ws.Range("I17").Value = 1.036780404994359
$ws.Range("J17").Value = 1.038497449443017
$ws.Range("K17").Value = 1.036525343384148
$ws.Range("L17").Value = 1.043746639471491
$ws.Range("M17").Value = 1.051293650830255
$ws.Range("N17").Value = 1.039972234938548
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032511827336232
$ws.Range("D18").Value = 1.033263025922973
$ws.Range("E18").Value = 1.040599869245524
$ws.Range("F18").Value = 1.048210129776615
$ws.Range("I18").Value = 1.036825806532653
$ws.Range("J18").Value = 1.038637715012515
$ws.Range("K18").Value = 1.036596400024568
$ws.Range("L18").Value = 1.043908061789959
$ws.Range("M18").Value = 1.051492611499242
$ws.Range("N18").Value = 1.040112699701244
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032587185192632
$ws.Range("D19").Value = 1.033301720729812
$ws.Range("E19").Value = 1.040669379963039
$ws.Range("F19").Value = 1.048292368679351
$ws.Range("I19").Value = 1.036841258845773
$ws.Range("J19").Value = 1.038685534541581
$ws.Range("K19").Value = 1.036620612557047
$ws.Range("L19").Value = 1.043963102885745
$ws.Range("M19").Value = 1.051560459398591
$ws.Range("N19").Value = 1.040160587139527
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032250185274275
$ws.Range("D20").Value = 1.033128663279929
$ws.Range("E20").Value = 1.040358563057358
$ws.Range("F20").Value = 1.047924653848582
$ws.Range("I20").Value = 1.036772040229242
$ws.Range("J20").Value = 1.038471645136183
$ws.Range("K20").Value = 1.036512265488498
$ws.Range("L20").Value = 1.043716947176911
$ws.Range("M20").Value = 1.051257057054113
$ws.Range("N20").Value = 1.039946393986639
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031154886527081
$ws.Range("D21").Value = 1.032565952934208
$ws.Range("E21").Value = 1.039348973031172
$ws.Range("F21").Value = 1.046730530460659
$ws.Range("I21").Value = 1.036545035917105
$ws.Range("J21").Value = 1.037775759717263
$ws.Range("K21").Value = 1.036158918274846
$ws.Range("L21").Value = 1.042916692033306
$ws.Range("M21").Value = 1.050271186276408
$ws.Range("N21").Value = 1.039249520330605
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030466476989015
$ws.Range("D22").Value = 1.032212098786499
$ws.Range("E22").Value = 1.038714902700798
$ws.Range("F22").Value = 1.045980778782357
$ws.Range("I22").Value = 1.036400804027548
$ws.Range("J22").Value = 1.037337841173921
$ws.Range("K22").Value = 1.035935916682976
$ws.Range("L22").Value = 1.042413559494899
$ws.Range("M22").Value = 1.049651732181176
$ws.Range("N22").Value = 1.038810979892701
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03083139221367
$ws.Range("D23").Value = 1.032399688310952
$ws.Range("E23").Value = 1.039050968819678
$ws.Range("F23").Value = 1.046378137393791
$ws.Range("I23").Value = 1.036477406560069
$ws.Range("J23").Value = 1.037570026616188
$ws.Range("K23").Value = 1.036054213152704
$ws.Range("L23").Value = 1.042680277476496
$ws.Range("M23").Value = 1.049980078171576
$ws.Range("N23").Value = 1.039043495064928
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03226855227598
$ws.Range("D24").Value = 1.0331380961262
$ws.Range("E24").Value = 1.04037550075065
$ws.Range("F24").Value = 1.047944691078471
$ws.Range("I24").Value = 1.036775820424293
$ws.Range("J24").Value = 1.03848330512626
$ws.Range("K24").Value = 1.036518175116621
$ws.Range("L24").Value = 1.043730363841696
$ws.Range("M24").Value = 1.051273592070237
$ws.Range("N24").Value = 1.03995807053524
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033937360946477
$ws.Range("D25").Value = 1.03399467910765
$ws.Range("E25").Value = 1.041915559309975
$ws.Range("F25").Value = 1.049767096949524
$ws.Range("I25").Value = 1.037115538158657
$ws.Range("J25").Value = 1.039541403527191
$ws.Range("K25").Value = 1.037052892831267
$ws.Range("L25").Value = 1.044948992351387
$ws.Range("M25").Value = 1.052776369008199
$ws.Range("N25").Value = 1.041017671557265
